# table9: drop the "Mexico 2010" / "Mexico 2020" data values (columns C & D)
# from every data row, and refresh the "Puerto Rico 2010" / "Puerto Rico 2020"
# (columns E & F) figures that changed with the recount.
#
# Values are written with a leading "'" so Excel stores them as text (matching
# the sheet's existing inline-string cells) instead of re-casting "0.49" etc.
# to a floating point number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $ws.Range($addr).Value = "'" + $value
}

# ---- Female ---------------------------------------------------------------

# Age
$ws.Range("C3:D3").ClearContents()
Set-TextValue "E3" "0.49"
Set-TextValue "F3" "0.49"

$ws.Range("C4:D4").ClearContents()
Set-TextValue "E4" "0.29"
Set-TextValue "F4" "0.29"

$ws.Range("C5:D5").ClearContents()
Set-TextValue "E5" "0.19"
Set-TextValue "F5" "0.18"

$ws.Range("C6:D6").ClearContents()
Set-TextValue "E6" "0.04"
Set-TextValue "F6" "0.05"

# Education Completed
$ws.Range("C8:D8").ClearContents()
Set-TextValue "E8" "0.27"

$ws.Range("C9:D9").ClearContents()
Set-TextValue "F9" "0.2"

$ws.Range("C10:D10").ClearContents()
Set-TextValue "E10" "0.35"

$ws.Range("C11:D11").ClearContents()
Set-TextValue "F11" "0.21"

# Household
$ws.Range("C13:D13").ClearContents()

$ws.Range("C14:D14").ClearContents()
Set-TextValue "E14" "0.28"

$ws.Range("C15:D15").ClearContents()
Set-TextValue "E15" "0.3"

$ws.Range("C16:D16").ClearContents()
Set-TextValue "E16" "0.41"

# ---- Male -------------------------------------------------------------------

# Age
$ws.Range("C19:D19").ClearContents()
Set-TextValue "E19" "0.55"
Set-TextValue "F19" "0.55"

$ws.Range("C20:D20").ClearContents()
Set-TextValue "E20" "0.29"
Set-TextValue "F20" "0.29"

$ws.Range("C21:D21").ClearContents()
Set-TextValue "E21" "0.14"
Set-TextValue "F21" "0.14"

$ws.Range("C22:D22").ClearContents()

# Education Completed
$ws.Range("C24:D24").ClearContents()
Set-TextValue "F24" "0.14"

$ws.Range("C25:D25").ClearContents()
Set-TextValue "F25" "0.23"

$ws.Range("C26:D26").ClearContents()
Set-TextValue "F26" "0.44"

$ws.Range("C27:D27").ClearContents()
Set-TextValue "F27" "0.19"

# Household
$ws.Range("C29:D29").ClearContents()
Set-TextValue "F29" "2.29"

$ws.Range("C30:D30").ClearContents()
Set-TextValue "E30" "0.19"

$ws.Range("C31:D31").ClearContents()

$ws.Range("C32:D32").ClearContents()
Set-TextValue "F32" "0.6"
